$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = -0.676580331437285
$ws.Range("D4").Value = 0.7736733172507795

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.7252941169329628
$ws.Range("D5").Value = -0.7826116191843241

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = -0.6280277759726285
$ws.Range("D8").Value = 0.6377326822774342

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0.689399516064055
$ws.Range("D9").Value = -0.6847479235958919
